$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16, shifting existing rows 16-38 down to 17-39.
$ws.Rows.Item(16).Insert()

# Fill the new row 16 with the record that was inserted (sharing the same
# constant columns A,B,C,E-J as the surrounding Nectarin records for this
# market/product).
$ws.Cells.Item(16, 1).Value = 1
$ws.Cells.Item(16, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(16, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(16, 4).Value = 44544
$ws.Cells.Item(16, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(16, 5).Value = 15
$ws.Cells.Item(16, 6).Value = "Fruta"
$ws.Cells.Item(16, 7).Value = 100103
$ws.Cells.Item(16, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(16, 9).Value = 100103006
$ws.Cells.Item(16, 10).Value = "Nectarín"
$ws.Cells.Item(16, 11).Value = "Artic Pride"
$ws.Cells.Item(16, 12).Value = "Segunda"
$ws.Cells.Item(16, 13).Value = 270
$ws.Cells.Item(16, 14).Value = 18000
$ws.Cells.Item(16, 15).Value = 20000
$ws.Cells.Item(16, 16).Value = 19000
$ws.Cells.Item(16, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(16, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(16, 19).Value = 1056
$ws.Cells.Item(16, 20).Value = 18
